$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1) and "全部类型" (sheetId 4) both received the same
# set of updated "想去人数" (F) / "最低票价" (G) values.

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet ---
$ws1.Range("G3").Value = 49
$ws1.Range("F4").Value = 203
$ws1.Range("F5").Value = 1114
$ws1.Range("F6").Value = 8320
$ws1.Range("F7").Value = 8320
$ws1.Range("F9").Value = 215
$ws1.Range("F10").Value = 6940
$ws1.Range("F11").Value = 177
$ws1.Range("F12").Value = 5083
$ws1.Range("F13").Value = 5556
$ws1.Range("F16").Value = 353
$ws1.Range("F18").Value = 317
$ws1.Range("F25").Value = 9329
$ws1.Range("F27").Value = 1711
$ws1.Range("F28").Value = 1032
$ws1.Range("F31").Value = 1900
$ws1.Range("F38").Value = 247
$ws1.Range("F39").Value = 1222
$ws1.Range("F41").Value = 4872
$ws1.Range("F42").Value = 382
$ws1.Range("F49").Value = 932
$ws1.Range("F50").Value = 1278
$ws1.Range("F51").Value = 47

# --- 全部类型 sheet ---
$ws4.Range("G3").Value = 49
$ws4.Range("F4").Value = 203
$ws4.Range("F6").Value = 1114
$ws4.Range("F7").Value = 8320
$ws4.Range("F9").Value = 215
$ws4.Range("F10").Value = 6940
$ws4.Range("F11").Value = 177
$ws4.Range("F14").Value = 5083
$ws4.Range("F15").Value = 5556
$ws4.Range("F18").Value = 353
$ws4.Range("F20").Value = 317
$ws4.Range("F25").Value = 9329
$ws4.Range("F27").Value = 1711
$ws4.Range("F28").Value = 1032
$ws4.Range("F31").Value = 1900
$ws4.Range("F38").Value = 247
$ws4.Range("F39").Value = 1222
$ws4.Range("F41").Value = 4872
$ws4.Range("F42").Value = 382
$ws4.Range("F49").Value = 932
$ws4.Range("F50").Value = 1278
$ws4.Range("F51").Value = 47
